$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (B:O) to the right
$ws.Columns("A:A").Insert()

# Set the new header cell value
$ws.Range("A1").Value = "Lehrer_ID"

# Update the active selection to A2
$ws.Range("A2").Select()
